$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two worker rows (16 and 17) and refresh the "Salario Basico"
# value for the worker now on row 16, as part of updating the EC database
# with the new estado de cuenta data.
$ws.Range("C16").Value = "1143393607"
$ws.Range("D16").Value = "JAIRO ALONSO QUINTANA BARRIOS"
$ws.Range("E16").Value = "2109"
$ws.Range("F16").Value = 42240
$ws.Range("G16").Value = 3566405

$ws.Range("C17").Value = "1088307874"
$ws.Range("D17").Value = "DANIEL GALVIZ RUIZ"
$ws.Range("E17").Value = "2403"
$ws.Range("F17").Value = 8040
$ws.Range("G17").Value = 6029743
